$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same style as the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data column, rows 2-15, all zero
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
